$d = $word.ActiveDocument

$pairs = @(
    @("50×34=", "22×73="),
    @("47×93=", "14×19="),
    @("52×23=", "13×49="),
    @("47×26=", "93×48="),
    @("63×82=", "23×22="),
    @("24×16=", "59×92="),
    @("94×30=", "53×88="),
    @("94×12=", "13×79="),
    @("74×58=", "63×98="),
    @("28×50=", "13×39="),
    @("80×37=", "91×29="),
    @("63×86=", "12×91="),
    @("54×21=", "15×35="),
    @("53×40=", "60×69="),
    @("54×80=", "25×22="),
    @("90×88=", "74×11="),
    @("14×34=", "88×30="),
    @("12×86=", "25×17="),
    @("18×72=", "49×67="),
    @("21×63=", "34×31="),
    @("43×45=", "88×69="),
    @("47×34=", "81×78="),
    @("97×37=", "33×45="),
    @("31×17=", "89×90="),
    @("33×21=", "16×48=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
